# modify LogManager and UI
#
# 1. Remove the stray "_GoBack" bookmark that sits by itself in an empty
#    paragraph (left over from a previous edit location).
# 2. Change the text "HIT.sql  " to "HIT20200115.sql  " (keeping the two
#    trailing spaces).
# 3. Re-create the "_GoBack" bookmark immediately after the run that now
#    contains "HIT20200115.sql  ", marking the new edit location - which is
#    exactly what Word itself does automatically when a document is edited.

$d = $word.ActiveDocument

# --- Step 1: drop the old "_GoBack" bookmark (wherever it currently is) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: update the text content ---
$body = $d.Content
$replaced = $body.Find.Execute("HIT.sql", $false, $false, $false, $false, $false, `
                                $true, 1, $false, "HIT20200115.sql", 2)

if (-not $replaced) {
    throw "Could not find 'HIT.sql' to replace"
}

# --- Step 3: place the "_GoBack" bookmark right after the edited run ---
# A temporary marker is appended after the full run text (including the
# trailing spaces) so the insertion point used for the bookmark is not the
# very last character position of the paragraph (that specific boundary is
# unreliable for Bookmarks.Add in this runtime); the marker is stripped
# again immediately after the bookmark has been created.
$marked = $d.Content
$markOk = $marked.Find.Execute("HIT20200115.sql  ", $false, $false, $false, $false, `
                                $false, $true, 1, $false, "HIT20200115.sql  ZZMARKERZZ", 2)

if (-not $markOk) {
    throw "Could not locate the updated text to anchor the bookmark"
}

$located = $d.Content
$locateOk = $located.Find.Execute("HIT20200115.sql  ")

if (-not $locateOk) {
    throw "Could not re-locate the updated text"
}

$anchor = $d.Range($located.End, $located.End)
$d.Bookmarks.Add("_GoBack", $anchor)
$goBack = $d.Bookmarks.Item("_GoBack")

$marker = $d.Range($goBack.End, $goBack.End + 11)
$marker.Text = ""

Write-Host "Updated text and re-anchored _GoBack bookmark."
